$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number (46060 = 2026-02-07).
# Update every data row (2 through 245) from 46060 to 46061 (2026-02-08).
for ($r = 2; $r -le 245; $r++) {
    $ws.Cells.Item($r, 3).Value = 46061
}
